$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Tipo licencias descripción") is reclassified from a measure to a dimension,
# gains a skos:Concept datatype and a new mapping-file reference row.
$ws.Range("D5").Copy($ws.Range("D6"))

$ws.Range("D3").Value = "iaest-dimension:tipo-licencias-descripcion"
$ws.Range("D4").Value = "dim"
$ws.Range("D5").Value = "skos:Concept"
$ws.Range("D6").Value = "mapping-tipo-licencias-descripcion.xlsx"
